$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title cell (A1) to new timestamp
$ws.Cells.Item(1,1).Value = "Datos actualizados a 20 de Marzo de 2020 a las 13:46"

# Country data rows: row, country, totals, new, active, recovered, critical, deathsToday, deaths
$data = @(
  @(4, "China", 80967, 39, 71150, 6569, 2136, 3, 3248),
  @(5, "Italia", 41035, 0, 4440, 33190, 2498, 0, 3405),
  @(6, "España", 19980, 1903, 1588, 17390, 939, 171, 1002),
  @(7, "Iran", 19644, 1237, 5979, 12232, 0, 149, 1433),
  @(8, "Alemania", 17372, 2052, 115, 17213, 2, 0, 44),
  @(9, "Estados Unidos", 14366, 577, 125, 14024, 64, 10, 217),
  @(10, "Francia", 10995, 0, 1295, 9328, 1122, 0, 372),
  @(11, "Corea del Sur", 8652, 87, 2233, 6325, 59, 3, 94),
  @(12, "Suiza", 4898, 676, 15, 4840, 0, 0, 43),
  @(13, "Reino Unido", 3269, 0, 65, 3060, 20, 0, 144),
  @(14, "Paises Bajos", 2460, 0, 2, 2382, 45, 0, 76),
  @(15, "Austria", 2333, 154, 9, 2318, 14, 0, 6),
  @(16, "Belgica", 2257, 462, 204, 2016, 164, 16, 37),
  @(17, "Noruega", 1835, 45, 1, 1827, 27, 0, 7),
  @(18, "Suecia", 1456, 17, 16, 1429, 21, 0, 11),
  @(19, "Dinamarca", 1226, 75, 1, 1216, 37, 3, 9),
  @(20, "Malasia", 1030, 130, 87, 941, 26, 0, 2),
  @(21, "Portugal", 1020, 234, 5, 1009, 26, 2, 6),
  @(22, "Japon", 963, 20, 215, 715, 50, 0, 33),
  @(23, "Australia", 876, 120, 46, 823, 2, 0, 7),
  @(24, "Canada", 873, 0, 11, 850, 1, 0, 12),
  @(25, "Chequia", 774, 80, 4, 770, 6, 0, 0),
  @(26, "Crucero", 712, 0, 527, 178, 14, 0, 7),
  @(27, "Israel", 705, 28, 15, 690, 10, 0, 0),
  @(28, "Brasil", 647, 7, 2, 638, 18, 0, 7),
  @(29, "Irlanda", 557, 0, 5, 549, 6, 0, 3),
  @(30, "Luxemburgo", 484, 149, 6, 474, 1, 0, 4),
  @(31, "Pakistan", 467, 13, 13, 451, 0, 1, 3),
  @(32, "Grecia", 464, 0, 19, 439, 16, 0, 6),
  @(33, "Catar", 460, 0, 10, 450, 6, 0, 0),
  @(34, "Finlandia", 430, 30, 10, 420, 2, 0, 0),
  @(35, "Islandia", 409, 79, 5, 404, 1, 0, 0),
  @(36, "Singapur", 385, 40, 124, 261, 14, 0, 0),
  @(37, "Polonia", 378, 23, 13, 359, 3, 1, 6),
  @(38, "Indonesia", 369, 60, 17, 320, 0, 7, 32),
  @(39, "Turquia", 359, 0, 0, 355, 0, 0, 4),
  @(40, "Chile", 342, 0, 0, 342, 6, 0, 0),
  @(41, "Tailandia", 322, 50, 42, 279, 1, 0, 1),
  @(42, "Eslovenia", 319, 0, 0, 318, 6, 0, 1),
  @(43, "Rumania", 308, 31, 31, 277, 11, 0, 0),
  @(44, "Barein", 284, 5, 110, 173, 4, 0, 1),
  @(45, "Estonia", 283, 16, 1, 282, 1, 0, 0),
  @(46, "Arabia Saudita", 274, 0, 8, 266, 0, 0, 0),
  @(47, "Ecuador", 260, 0, 1, 255, 2, 1, 4),
  @(48, "Egipto", 256, 0, 42, 207, 0, 0, 7),
  @(49, "Hong Kong", 256, 48, 98, 154, 4, 0, 4),
  @(50, "Peru", 234, 0, 1, 230, 7, 2, 3),
  @(51, "Filipinas", 230, 13, 8, 204, 1, 1, 18),
  @(52, "India", 223, 29, 23, 195, 0, 1, 5),
  @(53, "Sudafrica", 202, 52, 0, 202, 0, 0, 0),
  @(54, "Rusia", 199, 0, 9, 189, 0, 0, 1),
  @(55, "Irak", 192, 0, 49, 130, 0, 0, 13),
  @(56, "Mexico", 164, 46, 4, 159, 1, 0, 1),
  @(57, "Libano", 163, 6, 4, 155, 3, 0, 4),
  @(58, "Kuwait", 159, 11, 22, 137, 5, 0, 0),
  @(59, "San Marino", 144, 0, 4, 126, 12, 0, 14),
  @(60, "Emiratos Arabes Unidos", 140, 0, 31, 109, 2, 0, 0),
  @(61, "Panama", 137, 0, 1, 135, 7, 0, 1),
  @(62, "Armenia", 136, 14, 1, 135, 2, 0, 0),
  @(63, "Taiwan", 135, 27, 28, 105, 0, 1, 2),
  @(64, "Colombia", 128, 20, 1, 127, 0, 0, 0),
  @(65, "Argentina", 128, 0, 3, 122, 0, 0, 3),
  @(66, "Eslovaquia", 124, 0, 0, 124, 2, 0, 0),
  @(67, "Serbia", 118, 15, 2, 116, 4, 0, 0),
  @(68, "Croacia", 113, 3, 5, 107, 0, 0, 1),
  @(69, "Bulgaria", 112, 5, 1, 108, 0, 0, 3),
  @(70, "Letonia", 111, 25, 1, 110, 0, 0, 0),
  @(71, "Uruguay", 94, 15, 0, 94, 0, 0, 0),
  @(72, "Argelia", 90, 0, 32, 48, 0, 1, 10),
  @(73, "Costa Rica", 89, 2, 0, 87, 2, 1, 2),
  @(74, "Vietnam", 87, 2, 16, 71, 0, 0, 0),
  @(75, "Hungria", 85, 12, 7, 75, 6, 2, 3),
  @(76, "Islas Feroe", 80, 8, 3, 77, 0, 0, 0),
  @(77, "Brunei", 78, 5, 1, 77, 2, 0, 0),
  @(78, "Principado de Andorra", 75, 1, 1, 74, 2, 0, 0),
  @(79, "Albania", 70, 6, 0, 68, 2, 0, 2),
  @(80, "Sri Lanka", 70, 10, 3, 67, 0, 0, 0),
  @(81, "Jordania", 69, 0, 1, 68, 0, 0, 0),
  @(82, "Bosnia y Herzegovina", 69, 5, 2, 67, 0, 0, 0),
  @(83, "Bielorrusia", 69, 18, 15, 54, 0, 0, 0),
  @(84, "Republica de Chipre", 67, 0, 0, 67, 1, 0, 0),
  @(85, "Republica de Macedonia", 67, 17, 1, 66, 1, 0, 0),
  @(86, "Marruecos", 66, 3, 2, 61, 1, 1, 3),
  @(87, "Malta", 64, 11, 2, 62, 1, 0, 0),
  @(88, "Tunez", 54, 15, 1, 52, 2, 0, 1),
  @(89, "Kazajistan", 49, 5, 0, 49, 0, 0, 0),
  @(90, "Lituania", 49, 1, 1, 48, 1, 0, 0),
  @(91, "Moldavia", 49, 0, 1, 47, 3, 0, 1),
  @(92, "Oman", 48, 0, 13, 35, 0, 0, 0),
  @(93, "Estado de Palestina", 48, 1, 17, 31, 0, 0, 0),
  @(94, "Camboya", 47, 10, 1, 46, 0, 0, 0),
  @(95, "Guadalupe", 45, 12, 0, 45, 0, 0, 0),
  @(96, "Azerbaiyan", 44, 0, 7, 36, 0, 0, 1),
  @(97, "Georgia", 43, 3, 1, 42, 1, 0, 0),
  @(98, "Venezuela", 42, 0, 0, 42, 0, 0, 0),
  @(99, "Burkina Faso", 40, 7, 4, 35, 0, 0, 1),
  @(100, "Nueva Zelanda", 39, 11, 0, 39, 0, 0, 0),
  @(101, "Senegal", 38, 2, 2, 36, 0, 0, 0),
  @(102, "Republica Dominicana", 34, 0, 0, 32, 0, 0, 2),
  @(103, "Uzbekistan", 33, 10, 0, 33, 0, 0, 0),
  @(104, "Martinica", 32, 9, 0, 31, 7, 0, 1),
  @(105, "Reunion", 28, 0, 0, 28, 0, 0, 0),
  @(106, "Liechtenstein", 28, 0, 0, 28, 0, 0, 0),
  @(107, "Ucrania", 26, 0, 1, 22, 0, 0, 3),
  @(108, "Honduras", 24, 12, 0, 24, 0, 0, 0),
  @(109, "Afganistan", 24, 2, 1, 23, 0, 0, 0),
  @(110, "Camerun", 20, 7, 2, 18, 0, 0, 0),
  @(111, "Banglades", 20, 2, 3, 16, 1, 0, 1),
  @(112, "Consejo Danes para los Refugiados", 18, 4, 0, 18, 0, 0, 0),
  @(113, "Bolivia", 17, 2, 0, 17, 0, 0, 0),
  @(114, "Macao", 17, 0, 10, 7, 0, 0, 0),
  @(115, "Cuba", 16, 5, 0, 15, 0, 0, 1),
  @(116, "Jamaica", 16, 1, 2, 13, 0, 0, 1),
  @(117, "Ghana", 16, 5, 8, 8, 0, 0, 0),
  @(118, "Guayana Francesa", 15, 0, 0, 15, 0, 0, 0),
  @(119, "Guam", 14, 2, 0, 14, 0, 0, 0),
  @(120, "Montenegro", 13, 0, 0, 13, 0, 0, 0),
  @(121, "Paraguay", 13, 0, 0, 13, 1, 0, 0),
  @(122, "Maldivas", 13, 0, 0, 13, 0, 0, 0),
  @(123, "Mauricio", 12, 5, 0, 12, 0, 0, 0),
  @(124, "Nigeria", 12, 0, 1, 11, 0, 0, 0),
  @(125, "Monaco", 11, 1, 0, 11, 0, 0, 0),
  @(126, "Ruanda", 11, 0, 0, 11, 0, 0, 0),
  @(127, "Polinesia Francesa", 11, 5, 0, 11, 0, 0, 0),
  @(128, "Gibraltar", 10, 0, 2, 8, 0, 0, 0),
  @(129, "Trinidad yTobago", 9, 0, 0, 9, 0, 0, 0),
  @(130, "Togo", 9, 8, 0, 9, 0, 0, 0),
  @(131, "Etiopia", 9, 2, 0, 9, 0, 0, 0),
  @(132, "Guatemala", 9, 0, 0, 8, 0, 0, 1),
  @(133, "Costa de Marfil", 9, 0, 1, 8, 0, 0, 0),
  @(134, "Puerto Rico", 8, 2, 0, 8, 0, 0, 0),
  @(135, "Seychelles", 7, 1, 0, 7, 0, 0, 0),
  @(136, "Kenia", 7, 0, 0, 7, 0, 0, 0),
  @(137, "Guinea Ecuatorial", 6, 0, 0, 6, 0, 0, 0),
  @(138, "Mongolia", 6, 0, 0, 6, 0, 0, 0),
  @(139, "Kirguistan", 6, 3, 0, 6, 0, 0, 0),
  @(140, "Tanzania", 6, 0, 0, 6, 0, 0, 0),
  @(141, "Barbados", 5, 0, 0, 5, 0, 0, 0),
  @(142, "Guyana", 5, 0, 0, 4, 0, 0, 1),
  @(143, "Aruba", 5, 0, 1, 4, 0, 0, 0),
  @(144, "Mayotte", 4, 0, 0, 4, 0, 0, 0),
  @(145, "Surinam", 4, 3, 0, 4, 0, 0, 0),
  @(146, "Namibia", 3, 0, 0, 3, 0, 0, 0),
  @(147, "San Martin (Parte Francesa)", 3, 0, 0, 3, 0, 0, 0),
  @(148, "Congo", 3, 0, 0, 3, 0, 0, 0),
  @(149, "Bahamas", 3, 0, 0, 3, 0, 0, 0),
  @(150, "San Bartolome", 3, 0, 0, 3, 0, 0, 0),
  @(151, "Islas Virgenes de los Estados Unidos", 3, 0, 0, 3, 0, 0, 0),
  @(152, "Gabon", 3, 0, 0, 2, 0, 1, 1),
  @(153, "Islas Caimanes", 3, 0, 0, 2, 0, 0, 1),
  @(154, "Curazao", 3, 0, 0, 2, 0, 0, 1),
  @(155, "Republica de Africa Central", 2, 1, 0, 2, 0, 0, 0),
  @(156, "Benin", 2, 0, 0, 2, 0, 0, 0),
  @(157, "Liberia", 2, 0, 0, 2, 0, 0, 0),
  @(158, "Mauritania", 2, 0, 0, 2, 0, 0, 0),
  @(159, "Butan", 2, 1, 0, 2, 0, 0, 0),
  @(160, "Zambia", 2, 0, 0, 2, 0, 0, 0),
  @(161, "Groenlandia", 2, 0, 0, 2, 0, 0, 0),
  @(162, "Santa Lucia", 2, 0, 0, 2, 0, 0, 0),
  @(163, "Nueva Caledonia", 2, 0, 0, 2, 0, 0, 0),
  @(164, "Bermudas", 2, 0, 0, 2, 0, 0, 0),
  @(165, "Haiti", 2, 2, 0, 2, 0, 0, 0),
  @(166, "Sudan", 2, 0, 0, 1, 0, 0, 1),
  @(167, "Republica de Yibuti", 1, 0, 0, 1, 0, 0, 0),
  @(168, "Guinea", 1, 0, 0, 1, 0, 0, 0),
  @(169, "Suazilandia", 1, 0, 0, 1, 0, 0, 0),
  @(170, "Antigua y Barbuda", 1, 0, 0, 1, 0, 0, 0),
  @(171, "San Martin (Parte Holandesa)", 1, 0, 0, 1, 0, 0, 0),
  @(172, "Cabo Verde", 1, 1, 0, 1, 0, 0, 0),
  @(173, "El Salvador", 1, 0, 0, 1, 0, 0, 0),
  @(174, "Nicaragua", 1, 0, 0, 1, 0, 0, 0),
  @(175, "Republica del Chad", 1, 0, 0, 1, 0, 0, 0),
  @(176, "Fiyi", 1, 0, 0, 1, 0, 0, 0),
  @(177, "Somalia", 1, 0, 0, 1, 0, 0, 0),
  @(178, "Montserrat", 1, 0, 0, 1, 0, 0, 0),
  @(179, "Gambia", 1, 0, 0, 1, 0, 0, 0),
  @(180, "Niger", 1, 0, 0, 1, 0, 0, 0),
  @(181, "Isla de Man", 1, 0, 0, 1, 0, 0, 0),
  @(182, "Papua Nueva Guinea", 1, 1, 0, 1, 0, 0, 0),
  @(183, "Santa Sede", 1, 0, 0, 1, 0, 0, 0),
  @(184, "Angola", 1, 1, 0, 1, 0, 0, 0),
  @(185, "San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0),
  @(186, "Nepal", 1, 0, 1, 0, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
